# Apply the change described by the diff:
#  - A new "Run 50" column of results is inserted right before the existing
#    "Mean" column (which was column AZ). The "Mean" column shifts to BA and
#    its values are recalculated to include the new Run 50 results.
#  - The new "Run 50" column (AZ) gets header text "Run 50" and a constant
#    value of 22.64043143 for every data row.
#  - The (now shifted) "Mean" column (BA) gets the recalculated mean
#    20.08242935 for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at AZ (52nd column). This pushes the existing "Mean"
# column (and its formatting) from AZ to BA.
$ws.Columns.Item(52).Insert()

# New column header (the inserted cell already carries over the correct
# header style from the column insert operation above).
$ws.Cells.Item(1, 52).Value = "Run 50"

# Fill in the new "Run 50" values and the recalculated "Mean" values for each
# of the 13 data rows (rows 2 through 14).
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 52).Value = 22.64043143
    $ws.Cells.Item($r, 53).Value = 20.08242935
}

Write-Output "done"
